$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5953.875
$ws.Range("I62").Value = 4050
$ws.Range("K62").Value = 4050
$ws.Range("M62").Value = -3426
$ws.Range("H64").Value = 3412.5
$ws.Range("J64").Value = 2800
$ws.Range("L64").Value = 2800
$ws.Range("N64").Value = -3296
$ws.Range("H65").Value = 5953.875
$ws.Range("I65").Value = 4050
$ws.Range("K65").Value = 20250
$ws.Range("M65").Value = -17130
$ws.Range("H67").Value = 3412.5
$ws.Range("J67").Value = 2800
$ws.Range("L67").Value = 2800
$ws.Range("N67").Value = -4516
$ws.Range("H74").Value = 6114.5713
$ws.Range("J74").Value = 7933.3335
$ws.Range("L74").Value = 7933.3335
$ws.Range("N74").Value = -9805.333500000001
$ws.Range("H77").Value = 6114.5713
$ws.Range("J77").Value = 7933.3335
$ws.Range("L77").Value = 39666.6675
$ws.Range("N77").Value = -49026.6675
$ws.Range("H106").Value = 4870.5
$ws.Range("I106").Value = 1981
$ws.Range("K106").Value = 1981
$ws.Range("M106").Value = -1350
$ws.Range("H107").Value = 1481.8
$ws.Range("I107").Value = 1907.2727
$ws.Range("K107").Value = 1907.2727
$ws.Range("M107").Value = 12.72730000000001
$ws.Range("H112").Value = 28572680
$ws.Range("J112").Value = 1352.6774
$ws.Range("L112").Value = 4058.0322
$ws.Range("N112").Value = -6274.0322
$ws.Range("H113").Value = 3518.2856
$ws.Range("I113").Value = 1492.5
$ws.Range("K113").Value = 1492.5
$ws.Range("M113").Value = 1761.5
$ws.Range("H115").Value = 1588.421
$ws.Range("I115").Value = 1398.4615
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 4195.3845
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = -2628.3845
$ws.Range("N115").Value = -9134
$ws.Range("H137").Value = 1402733
$ws.Range("I137").Value = 2802850.2
$ws.Range("J137").Value = 2615.7058
$ws.Range("K137").Value = 8408550.600000001
$ws.Range("L137").Value = 7847.117400000001
$ws.Range("M137").Value = -8406000.600000001
$ws.Range("N137").Value = -12947.1174
$ws.Range("H138").Value = 4989.5454
$ws.Range("I138").Value = 953.56525
$ws.Range("J138").Value = 6210.9604
$ws.Range("K138").Value = 2860.69575
$ws.Range("L138").Value = 18632.8812
$ws.Range("M138").Value = 2279.30425
$ws.Range("N138").Value = -28912.8812
$ws.Range("H141").Value = 6351.857
$ws.Range("I141").Value = 7218.5312
$ws.Range("J141").Value = 3578.5
$ws.Range("K141").Value = 21655.5936
$ws.Range("L141").Value = 10735.5
$ws.Range("M141").Value = -16475.5936
$ws.Range("N141").Value = -21095.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1400.4584
$ws.Range("I110").Value = 1284.25
$ws.Range("J110").Value = 1981.5
$ws.Range("K110").Value = 1284.25
$ws.Range("L110").Value = 1981.5
$ws.Range("M110").Value = 760.75
$ws.Range("N110").Value = -6071.5
$ws.Range("H137").Value = 34536.668
$ws.Range("J137").Value = 45305
$ws.Range("L137").Value = 45305
$ws.Range("N137").Value = -55505

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 29800
$ws.Range("J55").Value = 29800
$ws.Range("L55").Value = 29800
$ws.Range("N55").Value = -30346
$ws.Range("H95").Value = 38410
$ws.Range("J95").Value = 38410
$ws.Range("L95").Value = 38410
$ws.Range("N95").Value = -43902
$ws.Range("H107").Value = 1422.1945
$ws.Range("I107").Value = 1335.5714
$ws.Range("J107").Value = 1725.375
$ws.Range("K107").Value = 1335.5714
$ws.Range("L107").Value = 1725.375
$ws.Range("M107").Value = 584.4286
$ws.Range("N107").Value = -5565.375
$ws.Range("H137").Value = 47750
$ws.Range("J137").Value = 47750
$ws.Range("L137").Value = 47750
$ws.Range("N137").Value = -57950

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3791.7
$ws.Range("I31").Value = 1391
$ws.Range("K31").Value = 1391
$ws.Range("M31").Value = -1096
$ws.Range("H34").Value = 3791.7
$ws.Range("I34").Value = 1391
$ws.Range("K34").Value = 1391
$ws.Range("M34").Value = -1189
$ws.Range("H122").Value = 3133.7778
$ws.Range("I122").Value = 1366.6666
$ws.Range("J122").Value = 4017.3333
$ws.Range("K122").Value = 4099.9998
$ws.Range("L122").Value = 12051.9999
$ws.Range("M122").Value = -1649.9998
$ws.Range("N122").Value = -16951.9999
$ws.Range("H141").Value = 31710.857
$ws.Range("J141").Value = 31710.857
$ws.Range("L141").Value = 31710.857
$ws.Range("N141").Value = -42070.857

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 2337.3333
$ws.Range("J82").Value = 2999.5
$ws.Range("L82").Value = 8998.5
$ws.Range("N82").Value = -9810.5
$ws.Range("H85").Value = 2337.3333
$ws.Range("J85").Value = 2999.5
$ws.Range("L85").Value = 8998.5
$ws.Range("N85").Value = -11806.5
$ws.Range("H113").Value = 652.1
$ws.Range("I113").Value = 602.125
$ws.Range("K113").Value = 1806.375
$ws.Range("M113").Value = 363.625
$ws.Range("H114").Value = 5562.5386
$ws.Range("I114").Value = 1828
$ws.Range("J114").Value = 5873.75
$ws.Range("K114").Value = 5484
$ws.Range("L114").Value = 17621.25
$ws.Range("M114").Value = -2230
$ws.Range("N114").Value = -24129.25
$ws.Range("H129").Value = 2668.1765
$ws.Range("I129").Value = 2295.2307
$ws.Range("J129").Value = 3880.25
$ws.Range("K129").Value = 6885.6921
$ws.Range("L129").Value = 11640.75
$ws.Range("M129").Value = -1885.6921
$ws.Range("N129").Value = -21640.75
$ws.Range("H134").Value = 3116.7778
$ws.Range("I134").Value = 2233.5
$ws.Range("J134").Value = 4883.3335
$ws.Range("K134").Value = 6700.5
$ws.Range("L134").Value = 14650.0005
$ws.Range("M134").Value = -1630.5
$ws.Range("N134").Value = -24790.0005
$ws.Range("H137").Value = 1530.1333
$ws.Range("J137").Value = 2430.5
$ws.Range("L137").Value = 7291.5
$ws.Range("N137").Value = -17491.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 31745.143
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 31745.143
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 31745.143
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -32057.143
$ws.Range("H113").Value = 1427.7941
$ws.Range("I113").Value = 1881.7142
$ws.Range("J113").Value = 1110.05
$ws.Range("K113").Value = 1881.7142
$ws.Range("L113").Value = 1110.05
$ws.Range("M113").Value = 288.2858000000001
$ws.Range("N113").Value = -5450.05
